$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1737.5
$ws.Range("I32").Value = 450
$ws.Range("K32").Value = 450
$ws.Range("M32").Value = -124
$ws.Range("H61").Value = 1635.5
$ws.Range("I61").Value = 1635.5
$ws.Range("K61").Value = 4906.5
$ws.Range("M61").Value = -4734.5
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H106").Value = 2727.3635
$ws.Range("I106").Value = 1123.625
$ws.Range("K106").Value = 1123.625
$ws.Range("M106").Value = -492.625
$ws.Range("H107").Value = 627.45
$ws.Range("I107").Value = 410.7143
$ws.Range("K107").Value = 410.7143
$ws.Range("M107").Value = 1509.2857
$ws.Range("H131").Value = 3146.3333
$ws.Range("J131").Value = 4429.6
$ws.Range("L131").Value = 13288.8
$ws.Range("N131").Value = -23368.8
$ws.Range("H132").Value = 876.37836
$ws.Range("I132").Value = 857.9143
$ws.Range("K132").Value = 2573.7429
$ws.Range("M132").Value = -43.74290000000019
$ws.Range("H135").Value = 633.58826
$ws.Range("I135").Value = 562.2857
$ws.Range("K135").Value = 5060.571300000001
$ws.Range("M135").Value = -2525.571300000001
$ws.Range("H138").Value = 3345.081
$ws.Range("I138").Value = 4525.467
$ws.Range("K138").Value = 13576.401
$ws.Range("M138").Value = -8436.400999999998

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3222.8872
$ws.Range("I32").Value = 2597.4917
$ws.Range("K32").Value = 2597.4917
$ws.Range("M32").Value = -2310.4917
$ws.Range("H61").Value = 2416.4092
$ws.Range("I61").Value = 896.5
$ws.Range("K61").Value = 896.5
$ws.Range("M61").Value = -684.5
$ws.Range("H74").Value = 1183.52
$ws.Range("I74").Value = 810.7059
$ws.Range("J74").Value = 1975.75
$ws.Range("K74").Value = 810.7059
$ws.Range("L74").Value = 1975.75
$ws.Range("M74").Value = 63.29409999999996
$ws.Range("N74").Value = -3723.75
$ws.Range("H77").Value = 1183.52
$ws.Range("I77").Value = 810.7059
$ws.Range("J77").Value = 1975.75
$ws.Range("K77").Value = 4053.5295
$ws.Range("L77").Value = 9878.75
$ws.Range("M77").Value = 314.4704999999999
$ws.Range("N77").Value = -18614.75
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()
$ws.Range("H136").Value = 2416.4092
$ws.Range("I136").Value = 896.5
$ws.Range("K136").Value = 2689.5
$ws.Range("M136").Value = -139.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2230.5454
$ws.Range("J86").Value = 2210
$ws.Range("L86").Value = 2210
$ws.Range("N86").Value = -4456
$ws.Range("H89").Value = 2230.5454
$ws.Range("J89").Value = 2210
$ws.Range("L89").Value = 11050
$ws.Range("N89").Value = -22282
$ws.Range("H94").Value = 928
$ws.Range("I94").Value = 810.8889
$ws.Range("K94").Value = 810.8889
$ws.Range("M94").Value = -359.8889
$ws.Range("H99").Value = 466
$ws.Range("I99").Value = 199.5
$ws.Range("K99").Value = 199.5
$ws.Range("M99").Value = 1298.5
$ws.Range("H105").Value = 2674.875
$ws.Range("I105").Value = 2699.8667
$ws.Range("K105").Value = 2699.8667
$ws.Range("M105").Value = -952.8667
$ws.Range("H134").Value = 8207.071
$ws.Range("I134").Value = 8542.23
$ws.Range("J134").Value = 3850
$ws.Range("K134").Value = 25626.69
$ws.Range("L134").Value = 11550
$ws.Range("M134").Value = -23091.69
$ws.Range("N134").Value = -16620

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1465.6
$ws.Range("I31").Value = 782.7368
$ws.Range("K31").Value = 782.7368
$ws.Range("M31").Value = -487.7368
$ws.Range("H34").Value = 1465.6
$ws.Range("I34").Value = 782.7368
$ws.Range("K34").Value = 782.7368
$ws.Range("M34").Value = -580.7368
$ws.Range("H58").Value = 2416804
$ws.Range("I58").Value = 3624462.2
$ws.Range("K58").Value = 3624462.2
$ws.Range("M58").Value = -3624259.2
$ws.Range("H62").Value = 5416.3
$ws.Range("I62").Value = 6827.8335
$ws.Range("J62").Value = 3299
$ws.Range("K62").Value = 6827.8335
$ws.Range("L62").Value = 3299
$ws.Range("M62").Value = -6203.8335
$ws.Range("N62").Value = -4547
$ws.Range("H65").Value = 5416.3
$ws.Range("I65").Value = 6827.8335
$ws.Range("J65").Value = 3299
$ws.Range("K65").Value = 34139.1675
$ws.Range("L65").Value = 16495
$ws.Range("M65").Value = -31019.1675
$ws.Range("N65").Value = -22735
$ws.Range("H136").Value = 2416804
$ws.Range("I136").Value = 3624462.2
$ws.Range("K136").Value = 10873386.6
$ws.Range("M136").Value = -10870836.6
$ws.Range("H141").Value = 65773.5
$ws.Range("J141").Value = 65773.5
$ws.Range("L141").Value = 65773.5
$ws.Range("N141").Value = -76133.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2949.75
$ws.Range("J68").Value = 2990.9565
$ws.Range("L68").Value = 8972.869499999999
$ws.Range("N68").Value = -10594.8695
$ws.Range("H71").Value = 2949.75
$ws.Range("J71").Value = 2990.9565
$ws.Range("L71").Value = 26918.6085
$ws.Range("N71").Value = -35030.6085
$ws.Range("H107").Value = 1416.871
$ws.Range("J107").Value = 1419
$ws.Range("L107").Value = 4257
$ws.Range("N107").Value = -8097
$ws.Range("H113").Value = 1150.6666
$ws.Range("I113").Value = 2935.5
$ws.Range("J113").Value = 640.7143
$ws.Range("K113").Value = 8806.5
$ws.Range("L113").Value = 1922.1429
$ws.Range("M113").Value = -6636.5
$ws.Range("N113").Value = -6262.1429
$ws.Range("H131").Value = 10432318
$ws.Range("J131").Value = 16665.178
$ws.Range("L131").Value = 49995.534
$ws.Range("N131").Value = -60075.534
$ws.Range("H132").Value = 1672.625
$ws.Range("J132").Value = 1753.3846
$ws.Range("L132").Value = 15780.4614
$ws.Range("N132").Value = -20840.4614
$ws.Range("H140").Value = 1868.8846
$ws.Range("I140").Value = 808.5217
$ws.Range("K140").Value = 2425.5651
$ws.Range("M140").Value = 2754.4349

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2884.3333
$ws.Range("I102").Value = 3065.4
$ws.Range("J102").Value = 2522.2
$ws.Range("K102").Value = 3065.4
$ws.Range("L102").Value = 2522.2
$ws.Range("M102").Value = -1443.4
$ws.Range("N102").Value = -5766.2
$ws.Range("H122").Value = 2402.5
$ws.Range("I122").Value = 1601
$ws.Range("J122").Value = 3043.7
$ws.Range("K122").Value = 4803
$ws.Range("L122").Value = 9131.099999999999
$ws.Range("M122").Value = -2353
$ws.Range("N122").Value = -14031.1
$ws.Range("H132").Value = 2027206.1
$ws.Range("I132").Value = 2565927.8
$ws.Range("K132").Value = 7697783.399999999
$ws.Range("M132").Value = -7695253.399999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1936.6923
$ws.Range("I68").Value = 1837.7
$ws.Range("K68").Value = 1837.7
$ws.Range("M68").Value = -1088.7
$ws.Range("H71").Value = 1936.6923
$ws.Range("I71").Value = 1837.7
$ws.Range("K71").Value = 9188.5
$ws.Range("M71").Value = -5444.5
$ws.Range("H81").Value = 49999
$ws.Range("J81").Value = 49999
$ws.Range("L81").Value = 49999
$ws.Range("N81").Value = -51995
$ws.Range("H84").Value = 49999
$ws.Range("J84").Value = 49999
$ws.Range("L84").Value = 149997
$ws.Range("N84").Value = -159981
$ws.Range("H122").Value = 6867
$ws.Range("I122").Value = 7397.1113
$ws.Range("K122").Value = 22191.3339
$ws.Range("M122").Value = -19741.3339
$ws.Range("H132").Value = 3301.487
$ws.Range("I132").Value = 1091.9445
$ws.Range("K132").Value = 3275.8335
$ws.Range("M132").Value = -745.8335000000002
$ws.Range("H135").Value = 45607
$ws.Range("J135").Value = 45607
$ws.Range("L135").Value = 45607
$ws.Range("N135").Value = -55747
$ws.Range("H136").Value = 2963.8223
$ws.Range("J136").Value = 4559.6
$ws.Range("L136").Value = 13678.8
$ws.Range("N136").Value = -18778.8

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 111960.7
$ws.Range("I122").Value = 126626.13
$ws.Range("J122").Value = 1970
$ws.Range("K122").Value = 379878.39
$ws.Range("L122").Value = 5910
$ws.Range("M122").Value = -377428.39
$ws.Range("N122").Value = -10810
$ws.Range("H126").Value = 4913.298
$ws.Range("I126").Value = 4876.5947
$ws.Range("J126").Value = 5049.1
$ws.Range("K126").Value = 14629.7841
$ws.Range("L126").Value = 15147.3
$ws.Range("M126").Value = -12159.7841
$ws.Range("N126").Value = -20087.3
